$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.02, 1.055303966889797, 1.060442637466814, 1.068286362931247, 1.074213382735155, 1.051295596365869, 1.06031152967956, 1.063169404769189, 1.07099195175525, 1.076903190060427, 1.061817293671138),
    @(1.02, 1.056254642602417, 1.06119572069401, 1.069197440842335, 1.075141284597221, 1.051575462635624, 1.060913824866073, 1.06373706922099, 1.071718741944794, 1.077647901902843, 1.0624204441859),
    @(1.02, 1.056870293256515, 1.0616834095119, 1.069787797989484, 1.075742555340395, 1.051755473929342, 1.061303398904246, 1.064104112766537, 1.072189216198378, 1.078130002388835, 1.062810571463891),
    @(1.02, 1.057129231427044, 1.061888526768955, 1.070036181179395, 1.075995533123257, 1.051830891304968, 1.061467138654798, 1.064258351440549, 1.07238704848887, 1.07833272941691, 1.06297454374367),
    @(1.02, 1.057172715211946, 1.061922972266039, 1.070077897324811, 1.076038021102159, 1.051843538986983, 1.061494629071215, 1.064284244867901, 1.072420268001434, 1.078366771190653, 1.063002073199631),
    @(1.02, 1.056873752736158, 1.061686149935477, 1.069791116125291, 1.075745934843816, 1.051756482679912, 1.061305586950093, 1.064106173973534, 1.072191859468077, 1.078132711034538, 1.062812762617015),
    @(1.02, 1.055625148063477, 1.060697063155417, 1.068594093537783, 1.074526793384692, 1.051390401925563, 1.060515108884232, 1.063361305617789, 1.071237533374413, 1.077154822203385, 1.062021161981631),
    @(1.02, 1.053428821732236, 1.058957240237449, 1.066491195284448, 1.072385130120537, 1.05073706916053, 1.059121070296057, 1.062046703343928, 1.069557413098891, 1.07543341108958, 1.060625143698762),
    @(1.02, 1.051967270103258, 1.057799505286335, 1.065093652936953, 1.070961888646383, 1.050296001534027, 1.058191018633595, 1.061168983126053, 1.068438425311847, 1.074287053307336, 1.059693771256261),
    @(1.02, 1.05133504687127, 1.057298719182042, 1.064489559738121, 1.070346702105673, 1.050103714235249, 1.057788143561727, 1.060788620034123, 1.067954164365735, 1.073790978732751, 1.05929032405558),
    @(1.02, 1.051100307857825, 1.057112784486417, 1.064265331989304, 1.070118359023621, 1.050032095100311, 1.057638475221189, 1.060647291460654, 1.06777432967354, 1.073606761653936, 1.059140443168832),
    @(1.02, 1.051150655775285, 1.057152664488281, 1.064313422369688, 1.070167331919784, 1.05004746646627, 1.057670580596219, 1.060677608930717, 1.06781290293507, 1.073646274692616, 1.059172594137177),
    @(1.02, 1.051315641290063, 1.057283348131435, 1.064471021743796, 1.070327823829757, 1.050097798151366, 1.057775772385722, 1.060776938678237, 1.067939298323378, 1.073775750323965, 1.059277935311087),
    @(1.02, 1.05141730717236, 1.057363877219221, 1.064568145076585, 1.070426730043213, 1.050128783331649, 1.057840581551769, 1.060838133108042, 1.068017180210485, 1.073855530809992, 1.059342836513583),
    @(1.02, 1.052009242205711, 1.057832751895573, 1.065133766888201, 1.07100273955729, 1.050308735622296, 1.058217752905229, 1.06119422023885, 1.06847056985929, 1.074319982686052, 1.059720543493627),
    @(1.02, 1.05238071896643, 1.058127004950043, 1.065488849059976, 1.071364346936537, 1.050421266720965, 1.058454301263437, 1.061417503452601, 1.068755041773758, 1.074611403786758, 1.059957427777636),
    @(1.02, 1.052597456595563, 1.058298687936436, 1.065696063883057, 1.071575371133775, 1.050486778535163, 1.058592260750593, 1.061547711334884, 1.068920995227391, 1.07478141422324, 1.060095583183088),
    @(1.02, 1.052671368931979, 1.058357235865596, 1.06576673594147, 1.071647342654202, 1.050509095037835, 1.058639298769751, 1.061592103835564, 1.068977585408996, 1.074839388347115, 1.060142688001629),
    @(1.02, 1.052340856655467, 1.058095429184436, 1.065450741605509, 1.071325539044248, 1.050409206193857, 1.058428923425509, 1.061393550314967, 1.068724517953106, 1.074580134007332, 1.059932013900267),
    @(1.02, 1.051267054487058, 1.057244862848845, 1.064424608246289, 1.070280558407919, 1.050082982105013, 1.057744796642035, 1.060747689771949, 1.067902076880586, 1.073737621668059, 1.059246915578291),
    @(1.02, 1.050592473233273, 1.056710539004475, 1.063780360159286, 1.069624490554983, 1.04987674373167, 1.057314528676516, 1.060341352973107, 1.067385215521279, 1.073208173830311, 1.058816036582909),
    @(1.02, 1.050950028043804, 1.056993750048274, 1.064121800354715, 1.069972193750764, 1.049986181362801, 1.057542633909106, 1.060556783995518, 1.06765919039202, 1.073488817875691, 1.059044465751093),
    @(1.02, 1.052358868527912, 1.058109696758907, 1.065467960410687, 1.07134307433684, 1.050414656215264, 1.058440390625038, 1.061404373797859, 1.068738310274161, 1.074594263383807, 1.059943497384535),
    @(1.02, 1.053996159319252, 1.059406653597585, 1.067034077097912, 1.072938008358619, 1.05090694555797, 1.059481588958292, 1.062386796295663, 1.069991576874961, 1.07587822202753, 1.060986174338849)
)

$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)  # B,C,D,E,F,I,J,K,L,M,N

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($row, $cols[$j]).Value = $rowData[$j]
    }
}

Write-Host "Updated $($data.Count) rows"